# Release 4.3.0: update to atlas-data 4.8.0
# - Refresh the x-mitre-tactic STIX IDs in column B (rows 2-15)
# - Refresh the "created"/"last modified" dates in columns F/G (rows 2-15)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tactics")

$newStixIds = @(
    "x-mitre-tactic--b08d83b4-fdc7-4e5f-93ef-f4479a2f2a14",
    "x-mitre-tactic--a467ab89-8c73-46bf-964a-9aa30822c85d",
    "x-mitre-tactic--a24a9db4-eb43-4017-95ae-36e4157cec6a",
    "x-mitre-tactic--46dee8aa-a49e-4177-8364-8cba270649e8",
    "x-mitre-tactic--4d802237-5f91-441b-bbd7-948b1eabbf6b",
    "x-mitre-tactic--deb1ca66-52d8-4ba9-b48d-eae1182e45b5",
    "x-mitre-tactic--72baacdb-3212-4cce-9785-9c068f2f3fd3",
    "x-mitre-tactic--c2ff12db-196b-4c2c-b05d-20f5423fe4f5",
    "x-mitre-tactic--5ff2c742-46ae-42ef-9d7b-0c1d76841197",
    "x-mitre-tactic--55db50f4-add2-4407-a20c-dc968ac3189d",
    "x-mitre-tactic--cb60f8f6-55cd-4ab3-a54e-76c4973c0e9e",
    "x-mitre-tactic--c920f814-50c3-4d9b-b0c9-8a910b29d357",
    "x-mitre-tactic--3fa0e3a6-3d30-4407-864d-7975b4feda3c",
    "x-mitre-tactic--56f79b50-5220-4170-bcbe-d07be1f8c8b5"
)

for ($i = 0; $i -lt $newStixIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newStixIds[$i]
}

$newDate = "14 March 2025"
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 6).Value = $newDate
    $ws.Cells.Item($row, 7).Value = $newDate
}
